$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "290.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.10%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-5.63%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.943"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.41%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07201"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-8.03%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.776"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-10.13%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.667"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.52%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.90%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8955"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.37%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1657"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.46%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07730"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.21%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08011"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-7.44%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03043"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.30%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1001"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.21%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001504"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.73%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005698"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.43%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.476"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.37%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.080"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.44%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3314"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.48%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1320"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.11%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.044"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-6.14%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2100"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.43%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04520"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.40%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004009"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-9.61%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.08%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01601"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.16%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04396"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-7.77%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007367"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.12%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1308"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.52%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007668"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002060"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-12.10%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009203"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-21.05%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005915"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-5.41%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.15%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "173.66%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002999"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.38%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.15%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.15%"
